$d = $word.ActiveDocument

# Locate the answer paragraph for the "Group 2, More Math" intersection problem
# (it contains the final answer "(6, 16)"); the target paragraph we need to
# fill in is the empty paragraph immediately following it.
$found = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "\(6,\s*16\)") {
        $found = $d.Paragraphs.Item($i + 1)
    }
}

if ($found -eq $null) {
    throw "Could not locate target paragraph"
}

$r = $found.Range

# Rebuild the (currently empty) paragraph with the paragraph-mark formatting
# (now red) plus a new run carrying the answer text, matching the
# surrounding Times New Roman / sz 24 / red styling used throughout the doc.
$openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' +
    '<w:p w14:paraId="7EAACCBF" w14:textId="2BC8ED24" w:rsidR="00FB719D" w:rsidRDefault="00FB719D" w:rsidP="000E2661">' +
    '<w:pPr>' +
    '<w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:color w:val="FF0000"/>' +
    '<w:sz w:val="24"/>' +
    '<w:szCs w:val="24"/>' +
    '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
    '<w:rPr>' +
    '<w:rFonts w:ascii="Times New Roman" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
    '<w:color w:val="FF0000"/>' +
    '<w:sz w:val="24"/>' +
    '<w:szCs w:val="24"/>' +
    '</w:rPr>' +
    '<w:t>Also intersect at (3,7) and (4,8), but those result in gibberish when decoded.</w:t>' +
    '</w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$r.InsertXML($openXml)

Write-Output "Updated paragraph with problem 2 bug fix text"
